$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2" = 93.62780766666667
    "H2" = 280.883423
    "I2" = 0.3228593149748609
    "J2" = 0.3228593149748609
    "M2" = 13.35941066666667
    "N2" = 40.078232
    "O2" = 0.4925555025958562
    "P2" = 0.4925555025958562
    "Q2" = 1250.812332438682
    "R2" = 11257.31099194814
    "S2" = 0.1590261321551965
    "T2" = 0.1590261321551964
    "G3" = 93.62780766666667
    "H3" = 280.883423
    "I3" = 0.3228593149748609
    "J3" = 0.3228593149748609
    "O3" = 0.03774352140193379
    "P3" = 0.03774352140193379
    "Q3" = 95.84719242886656
    "R3" = 862.624731859799
    "S3" = 0.01218584746456734
    "T3" = 0.01218584746456734
    "G4" = 93.62780766666667
    "H4" = 280.883423
    "I4" = 0.3228593149748609
    "J4" = 0.3228593149748609
    "M4" = 12.73953533333333
    "N4" = 38.218606
    "O4" = 0.4697009760022101
    "P4" = 0.46970097600221
    "Q4" = 1192.774763952038
    "R4" = 10734.97287556834
    "S4" = 0.1516473353550971
    "T4" = 0.1516473353550971
    "G5" = 66.39541
    "I5" = 0.228953097635189
    "J5" = 0.228953097635189
    "M5" = 13.35941066666667
    "N5" = 40.078232
    "O5" = 0.4925555025958562
    "P5" = 0.4925555025958562
    "Q5" = 887.0035485717067
    "R5" = 7983.03193714536
    "S5" = 0.1127721080765786
    "T5" = 0.1127721080765786
    "G6" = 66.39541
    "I6" = 0.228953097635189
    "J6" = 0.228953097635189
    "O6" = 0.03774352140193379
    "P6" = 0.03774352140193379
    "Q6" = 67.96926893044333
    "R6" = 611.7234203739899
    "S6" = 0.008641496140632792
    "T6" = 0.00864149614063279
    "G7" = 66.39541
    "I7" = 0.228953097635189
    "J7" = 0.228953097635189
    "M7" = 12.73953533333333
    "N7" = 38.218606
    "O7" = 0.4697009760022101
    "P7" = 0.46970097600221
    "Q7" = 845.8466716661534
    "R7" = 7612.62004499538
    "S7" = 0.1075394934179775
    "T7" = 0.1075394934179775
    "G8" = 129.9724656666667
    "H8" = 389.917397
    "I8" = 0.4481875873899502
    "J8" = 0.4481875873899502
    "M8" = 13.35941066666667
    "N8" = 40.078232
    "O8" = 0.4925555025958562
    "P8" = 0.4925555025958562
    "Q8" = 1736.355544200234
    "R8" = 15627.1998978021
    "S8" = 0.2207572623640812
    "T8" = 0.2207572623640812
    "G9" = 129.9724656666667
    "H9" = 389.917397
    "I9" = 0.4481875873899502
    "J9" = 0.4481875873899502
    "O9" = 0.03774352140193379
    "P9" = 0.03774352140193379
    "Q9" = 133.0533763169846
    "R9" = 1197.480386852861
    "S9" = 0.01691617779673366
    "T9" = 0.01691617779673366
    "G10" = 129.9724656666667
    "H10" = 389.917397
    "I10" = 0.4481875873899502
    "J10" = 0.4481875873899502
    "M10" = 12.73953533333333
    "N10" = 38.218606
    "O10" = 0.4697009760022101
    "P10" = 0.46970097600221
    "Q10" = 1655.788818720954
    "R10" = 14902.09936848858
    "S10" = 0.2105141472291354
    "T10" = 0.2105141472291354
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
